$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "36.337.81"
$ws.Range("E2").Value = "  -0.58%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.019.29"
$ws.Range("E3").Value = "  +2.18%  "

$ws.Range("E4").Value = "  -0.13%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "230.58"
$ws.Range("E5").Value = "  -8.39%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.597"
$ws.Range("E6").Value = "  -1.92%  "

$ws.Range("E7").Value = "  +0.05%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "54.30"
$ws.Range("E8").Value = "  -1.44%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.368"
$ws.Range("E9").Value = "  -1.56%  "

$ws.Range("E10").Value = "  +3.11%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0744"
$ws.Range("E11").Value = "  -2.13%  "

$ws.Range("E12").Value = "  -0.89%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.321.67"
$ws.Range("E13").Value = "  +1.97%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "14.15"
$ws.Range("E14").Value = "  -0.31%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "19.95"
$ws.Range("E15").Value = "  -6.96%  "

$ws.Range("E16").Value = "  -3.32%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.08"
$ws.Range("E17").Value = "  -1.28%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.993.38"
$ws.Range("E18").Value = "  +0.38%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "36.541.30"
$ws.Range("E19").Value = "  +0.15%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "67.34"
$ws.Range("E20").Value = "  -3.54%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.44"
$ws.Range("E21").Value = "  +8.55%  "

$ws.Range("E22").Value = "  -3.21%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "219.82"
$ws.Range("E23").Value = "  -5.76%  "

$ws.Range("E24").Value = "  +0.04%  "

$ws.Range("E25").Value = "  +0.40%  "

$ws.Range("E26").Value = "  -6.67%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "162.66"
$ws.Range("E27").Value = "  -0.52%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.58"
$ws.Range("E28").Value = "  -2.21%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.127"
$ws.Range("E29").Value = "  +3.27%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "18.74"
$ws.Range("E30").Value = "  -2.56%  "

$ws.Range("E31").Value = "  +0.93%  "

$ws.Range("E32").Value = "  -0.96%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.34"
$ws.Range("E33").Value = "  -3.40%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0597"
$ws.Range("E34").Value = "  -4.96%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.44"
$ws.Range("E35").Value = "  +4.20%  "

$ws.Range("E36").Value = "  -2.24%  "

$ws.Range("E37").Value = "  +0.15%  "

$ws.Range("E38").Value = "  -2.82%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.25"
$ws.Range("E39").Value = "  -5.84%  "

$ws.Range("E40").Value = "  +5.65%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.91"
$ws.Range("E41").Value = "  -2.89%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.460.49"
$ws.Range("E42").Value = "  +2.44%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0928"
$ws.Range("E43").Value = "  +2.35%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "92.55"
$ws.Range("E44").Value = "  +5.21%  "

$ws.Range("E45").Value = "  -2.14%  "

$ws.Range("E46").Value = "  -4.88%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "15.39"
$ws.Range("E47").Value = "  -0.69%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.94"
$ws.Range("E48").Value = "  +30.61%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.994"
$ws.Range("E49").Value = "  -1.51%  "

$ws.Range("E50").Value = "  +0.19%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.82"
$ws.Range("E51").Value = "  -0.35%  "
